# Apply the "linear equations in newborn outcomes and pregnancy supervisor" edit:
#  - Add two new parameter rows to the parameter_values sheet:
#      prob_enceph_severity      (inserted after prob_encephalopathy, old row 8)
#      prob_retinopathy_severity (inserted after prob_retinopathy_preterm, old row 12)
#  - Apply the right-aligned "value" style (same as the existing B6 cell) to the
#    value cells of the rows that now sit just above/around the new severity rows
#    (B5, B7, B8, and the new B9).
#  - Update the sheet selection to reflect the area that was last worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row after the current row 8 (prob_encephalopathy) ---------
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "prob_enceph_severity"
$ws.Range("B9").Value = "[0.422, 0.338, 0.24]"
# No source cell for this new row.
$ws.Range("D9").Clear()

# --- Insert a new row after the (now shifted) row 13 (prob_retinopathy_preterm) ---
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "prob_retinopathy_severity"
$ws.Range("B14").Value = "[0.4, 0.3, 0.2, 0.1]"
# Same (dummy) source as prob_retinopathy_preterm directly above it.
$ws.Range("D13").Copy($ws.Range("D14"))

# --- Apply the right-aligned numeric style (as already used on B6) to the ---
# --- value cells around the newly-inserted row, matching the commit's     ---
# --- formatting cleanup.                                                  ---
$ws.Range("B6").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Reflect the last worked-on selection -----------------------------------
$ws.Range("D13:D14").Select()
